$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix capitalization / wording in existing Case 1 block ---
$ws.Range("B4").Value = "ATR_Amplitude = 3.5V"

# Remove the old "ARP = 250ms" row (row 11) from Case 1 block
$ws.Range("B11").ClearContents()

# --- Case 2 block edits ---
$ws.Range("C13").Value = "NAT_HEART decreases to 60"
$ws.Range("B15").Value = "VENT_Amplitude = 3.8V"
$ws.Range("B18").Value = "High Threshold"

# Remove the old "VRP = 250ms" row (row 22) at the end of Case 2 block
$ws.Range("B22").ClearContents()

# --- New Case 3 block (AAIR) rows 24-34 ---
$ws.Range("A24").Value = "Case 3"
$ws.Range("A24").Font.Bold = $true
$ws.Range("B24").Value = "Mode = 7 (AAIR)"
$ws.Range("C24").Value = "NAT_HEART increases to 100"

$ws.Range("B25").Value = "RATE = 100"
$ws.Range("B26").Value = "ATR_Amplitude = 3V"
$ws.Range("B27").Value = "ATR_PULSE_WIDTH = 8ms"
$ws.Range("B28").Value = "NAT_HEART = 90"
$ws.Range("B29").Value = "Medium Threshold"
$ws.Range("B30").Value = "Response_Factor = 8"
$ws.Range("B31").Value = "Response_time = 10sec"
$ws.Range("B32").Value = "Recovery_time = 10sec"
$ws.Range("B33").Value = "ARP = 250ms"
$ws.Range("B34").Value = "ATRSENSITIVITY_Amplitude = 1.5mV "

# --- New Case 4 block (VVIR) rows 35-46 ---
$ws.Range("A35").Value = "Case 4"
$ws.Range("A35").Font.Bold = $true
$ws.Range("B36").Value = "Mode 8 = VVIR"
$ws.Range("C36").Value = "Boundry case, NAT_HEART remains constant at 65"

$ws.Range("B37").Value = "RATE = 65"
$ws.Range("B38").Value = "VENT_Amplitude = 5V"
$ws.Range("B39").Value = "VENT_PULSE_WIDTH = 10ms"
$ws.Range("B40").Value = "NAT_HEART = 65"
$ws.Range("B41").Value = "High Threshold"
$ws.Range("B42").Value = "Response_Factor = 5"
$ws.Range("B43").Value = "Response_time = 5sec"
$ws.Range("B44").Value = "Recovery_time = 5sec"
$ws.Range("B45").Value = "VRP = 300"
$ws.Range("B46").Value = "VENTSENSITIVITY_Amplitude = 2mV "

# --- Update selection to match final state (active cell A24) ---
$ws.Range("A24").Select()
